# Rewrite the 'Export' sheet's data region to match the refreshed Saldo export:
# new accounts inserted near the top, several zero/duplicate balances removed,
# two balances updated, and the whole list re-sorted by balance (desc).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture how many rows are currently used so we can wipe everything
# (old data rows + the old blank separator + the old footer row) before
# rewriting it fresh below. A small safety pad covers any stray cells
# just past the detected used range.
$oldLastRow = $ws.UsedRange.Rows.Count + 5
if ($oldLastRow -lt 5) { $oldLastRow = 5 }
$ws.Range("A2:C" + $oldLastRow).ClearContents()

# Conta, Nome, Saldo — ordered exactly as in the refreshed export
$data = @(
    @("005642649", "VR", 500000),
    @("005103059", "WALQUIRIA", 228900),
    @("005608744", "DOUGLAS", 175742.29),
    @("004352384", "BRASFORT", 152091.17),
    @("004974089", "CELIA", 118460.07),
    @("005277762", "NIVALDO", 114450),
    @("004361159", "HFR", 59491.2),
    @("004321016", "JOAQUIM", 31489.64),
    @("004895776", "FERNANDO", 30325.95),
    @("005064129", "THIAGO", 26552.57),
    @("004224011", "THOMAS", 16134.66),
    @("004364200", "BLOCO", 13305.55),
    @("004690692", "PHYLIA", 9906.32),
    @("004971783", "ERIKA", 5000),
    @("001761119", "BLUEMETRIX", 2016.95),
    @("005599726", "JORGE", 1960.29),
    @("004467884", "ANA", 1748.17),
    @("004854514", "MARCIA", 1200.84),
    @("004487140", "VALMIR", 1160),
    @("004392159", "RODRIGO", 900.21),
    @("004369172", "LUIZA", 679.85),
    @("005171652", "BRUNO", 506.45),
    @("004350197", "GISELA", 387.94),
    @("005135532", "FELIPE", 280.13),
    @("004498637", "TIAGO", 243.4),
    @("004958578", "ASSAKO", 110.07),
    @("005270025", "DENIZE", 100),
    @("004342617", "JURACI", 99.96),
    @("005173958", "VENIA", 99.86),
    @("004472076", "RUBENS", 99.18),
    @("004339183", "JALISON", 95.69),
    @("004517506", "LUIZ", 94.49),
    @("004643880", "GABRIEL", 94.09),
    @("004431591", "MARIO", 93.87),
    @("004027477", "GABRIELA", 93.28),
    @("005055226", "DANILO", 89.16),
    @("004278033", "DAISY", 86.72),
    @("004432579", "ANA", 86.5),
    @("004212132", "JOAO", 86.38),
    @("004207278", "CESAR", 84.93),
    @("004239387", "LUIZ", 83.93),
    @("003115072", "VICTOR", 81.27),
    @("005348011", "TATIANA", 80.91),
    @("005009947", "VERANICE", 80.79),
    @("004318604", "RENAN", 80.51),
    @("004994036", "BALTASAR", 80.5),
    @("000330949", "RENATO", 80.08),
    @("004809902", "PEDRO", 79.88),
    @("004267976", "E3", 79.84),
    @("005032151", "ANA", 79.27),
    @("004454365", "RAFAEL", 79.25),
    @("005256849", "SANDRO", 77.17),
    @("004479734", "RODRIGO", 76),
    @("004453045", "JULIAN", 75.83),
    @("004230529", "LAIS", 75.09),
    @("004340984", "RENATA", 73.59),
    @("004277637", "LARA", 73.51),
    @("004207374", "ANGELICA", 72.95),
    @("004520100", "ALEXANDRE", 71.72),
    @("004268684", "PATRICIA", 68.11),
    @("004332207", "IRACY", 67.97),
    @("004452507", "DANIELA", 65.08),
    @("004756968", "DANIELY", 64.89),
    @("004855596", "MARIANA", 64.36),
    @("005186167", "ANDREA", 63.64),
    @("005305965", "SIDMAR", 62.42),
    @("004335251", "EDMUNDO", 62.39),
    @("004242237", "MARIAH", 60.32),
    @("005348975", "JULIA", 60),
    @("004321092", "DANIEL", 58.87),
    @("004451996", "ADRIANO", 58.35),
    @("004381194", "ALINNE", 58.16),
    @("004472760", "SANDRA", 57.27),
    @("004804036", "LUCIANA", 56.61),
    @("004693308", "LAURA", 56.37),
    @("004459461", "INTERLAGOS", 56.02),
    @("004215217", "CAROLINA", 55.66),
    @("004575632", "ADELE", 54.96),
    @("005018038", "ELAINE", 54.42),
    @("005216881", "RENAN", 53.5),
    @("003836362", "ISABELLA", 51.76),
    @("004400640", "FELIPE", 51.44),
    @("004208447", "LEILA", 50),
    @("005077648", "DUNAS", 49.87),
    @("004229526", "EDUARDA", 49.12),
    @("004207184", "CRISTINA", 48.93),
    @("004466221", "WALTER", 48.66),
    @("004517080", "TATIANA", 47.35),
    @("004208733", "REINALDO", 46.33),
    @("004115403", "HEBERT", 45.24),
    @("001731007", "GUILHERME", 44.59),
    @("004272426", "RODRIGO", 44.35),
    @("005070742", "JUSCELINO", 44.06),
    @("005000460", "MARIANA", 41.03),
    @("004397124", "MURYLO", 40.43),
    @("004752615", "LUZIMAR", 39.1),
    @("004398174", "DANIELE", 39.08),
    @("004238164", "DANIELA", 38.3),
    @("001294033", "VIVIANE", 38.26),
    @("002401479", "JULIO", 37.84),
    @("005055239", "NORMAN", 37.66),
    @("004259649", "BENTO", 37.61),
    @("004265173", "JULIA", 36.93),
    @("004973881", "ISABELLA", 36.87),
    @("004329229", "GABRIEL", 36.25),
    @("004431689", "LUCA", 35.72),
    @("004207641", "MAGALI", 35.14),
    @("004243043", "SUELI", 34.15),
    @("005437764", "EVA", 33.57),
    @("004470679", "RODOLFO", 33.54),
    @("005035754", "JOSE", 33.36),
    @("005558076", "ALEXANDRE", 28.84),
    @("005274028", "RAFAEL", 28.62),
    @("004999434", "EDUARDO", 28.06),
    @("005366255", "RAPHAELA", 26.47),
    @("004377415", "ANGELA", 26.37),
    @("005591536", "GUSTAVO", 25.8),
    @("004940560", "CRISTIANO", 25.34),
    @("004240292", "MARCO", 24.3),
    @("004404724", "LEANDRO", 24.14),
    @("005581299", "ZILDA", 23.99),
    @("005009922", "ANA", 23.62),
    @("004259659", "BENTO", 22.99),
    @("004563237", "FERNANDO", 21.74),
    @("005219257", "CAROLINE", 21.37),
    @("004214604", "MARIA", 20.75),
    @("004204255", "AMADO", 18.77),
    @("004527606", "MARCIA", 18.68),
    @("004181486", "ANDREA", 18),
    @("004479463", "HENRIQUE", 17.36),
    @("004946997", "EDUARDO", 16.98),
    @("005143579", "GABRIEL", 16.18),
    @("005169333", "EDUARDO", 16.12),
    @("004920447", "MARILIA", 16.02),
    @("004995535", "ASIEL", 15.76),
    @("005133039", "PAULO", 15.48),
    @("004848279", "LEONARDO", 15.25),
    @("005135281", "RAFAEL", 15.2),
    @("004581652", "CINCO", 14.94),
    @("004422594", "WANDIR", 14.67),
    @("004472538", "RODOLFO", 14.25),
    @("002064834", "RAFAELA", 13.24),
    @("004207955", "SILVANIA", 12.48),
    @("005366671", "TATIANA", 11.45),
    @("004752461", "SERGIO", 10.77),
    @("004216298", "FLORDELIZ", 9.76),
    @("002894447", "JOAO", 9.17),
    @("004264780", "MARCELO", 8.99),
    @("004479965", "DIEGO", 8.78),
    @("004693631", "NELY", 7.36),
    @("004767746", "ISABELE", 7.33),
    @("004419765", "WALTER", 7.21),
    @("004508159", "FELIPE", 7.11),
    @("004530494", "ROSANGELA", 6.94),
    @("004854496", "JOSE", 6.64),
    @("004289402", "LARISSA", 6.47),
    @("004805133", "PATRICIA", 6.22),
    @("004448501", "JOAO", 5.55),
    @("005142624", "RODRIGO", 4.75),
    @("004462543", "RODOLFO", 4.3),
    @("005242683", "LUCAS", 4.26),
    @("004848927", "ULDARICO", 3.62),
    @("004936634", "LEONARDO", 3.08),
    @("004308815", "ZELI", 1.25),
    @("004431546", "GABRIELA", 0.97),
    @("004223502", "BRUNA", 0.78),
    @("005624274", "CLAYTON", 0.68),
    @("002694089", "VITOR", 0.65),
    @("004453302", "ISABELLA", 0.39),
    @("004806286", "VERA", 0.19),
    @("004371857", "NAZARETH", 0.18),
    @("004357159", "JOAO", 0.15),
    @("004320840", "NATALIA", 0.14),
    @("004335031", "EDMUNDO", 0.11),
    @("004466350", "RAQUEL", 0.11),
    @("005047946", "GABRIEL", 0.09),
    @("004589311", "CLARICE", 0.06),
    @("001090818", "MARIA", 0.01),
    @("002878817", "GUILHERME", 0.01),
    @("004384258", "PAULA", 0.01),
    @("004400000", "VILMA", 0.01),
    @("004474935", "MELISSA", 0.01),
    @("004475136", "MATHEO", 0.01),
    @("004612043", "YURI", 0.01),
    @("005314853", "GLEUBER", 0.01)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = "'" + $r[0]   # keep leading zeros as text
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $row++
}

# Row after the data stays blank (mirrors the source export layout),
# then the applied-filters footer note goes on the next row.
$footerRow = $row + 1
$ws.Cells.Item($footerRow, 1).Value = "Filtros aplicados:
DataFim é (Em branco)
nr_saldo_disponivel não é 0
Posição é Posição D-1
DataFim é (Em branco)
CARTEIRA não está em branco
NR_CONTA não está em branco
TIPO_LANCAMENTO não é ED, ET ou Liquidação Doador"

Write-Output ("New used rows: " + $ws.UsedRange.Rows.Count)
